# Auto-generated edit script: updates cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be
# auto-coerced to a number (losing formatting like trailing zeros).
$textCells = @("D5","D6","D10","D12","D14","D18","D21","D22","D23","D24","D25","D26","D27","D32","D37","D38","D42","D43","D44","D45","D46","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.898.56"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.907.31"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "591.14"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "144.95"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "2.905.14"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "33.50"
$ws.Range("E14").Value = "  -4.84%  "
$ws.Range("D16").Value = "3.389.92"
$ws.Range("D17").Value = "60.827.81"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -4.18%  "
$ws.Range("D19").Value = "2.902.00"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").Value = "13.52"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "7.08"
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("D24").Value = "81.30"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "10.90"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "11.99"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "7.07"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").Value = "0.0₃0853"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "5.60"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "8.61"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "0.290"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").Value = "40.12"
$ws.Range("E44").Value = "  -9.15%  "
$ws.Range("D45").Value = "375.23"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").Value = "0.0347"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("D47").Value = "2.707.02"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "129.56"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("D50").Value = "24.02"
$ws.Range("E50").Value = "  -8.27%  "
$ws.Range("E51").Value = "  -1.63%  "
